# Apply the docxtotei / teitodocx round-tripping style fixes:
#  - introduce four new custom paragraph styles (teisigned, teispeech,
#    GeneratedTitle, GeneratedSubTitle)
#  - re-point the (empty) first paragraph at the new GeneratedSubTitle
#    style instead of its ad-hoc hanging indent

$d = $word.ActiveDocument

# --- tei_signed ------------------------------------------------------
$teiSigned = $d.Styles.Add("teisigned", 1)
$teiSigned.NameLocal = "tei_signed"
$teiSigned.BaseStyle = $d.Styles("Normal")
$teiSigned.QuickStyle = $true
$teiSigned.ParagraphFormat.SpaceBefore = 18
$teiSigned.ParagraphFormat.LeftIndent = 21.55
$teiSigned.ParagraphFormat.FirstLineIndent = -21.55

# --- tei_speech --------------------------------------------------------
$teiSpeech = $d.Styles.Add("teispeech", 1)
$teiSpeech.NameLocal = "tei_speech"
$teiSpeech.BaseStyle = $d.Styles("Normal")
$teiSpeech.QuickStyle = $true
$teiSpeech.ParagraphFormat.LeftIndent = 21.6
$teiSpeech.ParagraphFormat.FirstLineIndent = -21.6

# --- GeneratedTitle ------------------------------------------------------
$generatedTitle = $d.Styles.Add("GeneratedTitle", 1)
$generatedTitle.BaseStyle = $d.Styles("Title")
$generatedTitle.QuickStyle = $true

# --- GeneratedSubTitle ----------------------------------------------------
$generatedSubTitle = $d.Styles.Add("GeneratedSubTitle", 1)
$generatedSubTitle.BaseStyle = $d.Styles("Subtitle")
$generatedSubTitle.QuickStyle = $true

# Re-point the lone (empty) paragraph at the new GeneratedSubTitle style,
# replacing its explicit hanging-indent pPr.
$p = $d.Paragraphs(1)
$p.Style = $d.Styles("GeneratedSubTitle")
